$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title text in B2
$rng = $ws.Range("B2")
$rng.Value = "DevOps zero to hero course notes"

# Bold, 16pt font for the title
$rng.Font.Bold = $true
$rng.Font.Size = 16

# Row grows to fit the larger title font
$ws.Rows.Item(2).RowHeight = 21

# Leave the selection where it was left in the authored file
$ws.Range("H15").Select() | Out-Null
